# Append 4 new data rows (230-233) to Sheet1, mirroring the existing
# row layout: col A = sequential index (bold/centered/bordered style,
# copied from the prior row), col B = date as literal text, cols C:W = numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds dates as plain text (not real Excel dates) in this workbook.
# Pre-format the target cells as Text so the "dd/mm/yyyy" strings are not
# auto-parsed into date serials, then strip the format again so the saved
# cells carry no extra number-format styling (matching the source rows).
$ws.Range("B230:B233").NumberFormat = "@"

# Row 230
$ws.Cells.Item(230, 1).Value = 228
$ws.Cells.Item(230, 2).Value = "30/08/2024"
$ws.Cells.Item(230, 3).Value = 198.8278
$ws.Cells.Item(230, 4).Value = 166.375
$ws.Cells.Item(230, 5).Value = 156.1
$ws.Cells.Item(230, 6).Value = 148
$ws.Cells.Item(230, 7).Value = 64.75
$ws.Cells.Item(230, 8).Value = 134
$ws.Cells.Item(230, 9).Value = 67.83799999999999
$ws.Cells.Item(230, 10).Value = 62
$ws.Cells.Item(230, 11).Value = 125
$ws.Cells.Item(230, 12).Value = 299.4042307692308
$ws.Cells.Item(230, 13).Value = 120
$ws.Cells.Item(230, 14).Value = 235
$ws.Cells.Item(230, 15).Value = 198
$ws.Cells.Item(230, 16).Value = 181
$ws.Cells.Item(230, 17).Value = 175
$ws.Cells.Item(230, 18).Value = 99
$ws.Cells.Item(230, 19).Value = 160
$ws.Cells.Item(230, 20).Value = 0.3192307692307693
$ws.Cells.Item(230, 21).Value = 64.48999999999999
$ws.Cells.Item(230, 22).Value = 118.25
$ws.Cells.Item(230, 23).Value = 64.48999999999999

# Row 231
$ws.Cells.Item(231, 1).Value = 229
$ws.Cells.Item(231, 2).Value = "02/09/2024"
$ws.Cells.Item(231, 3).Value = 204.4669696969697
$ws.Cells.Item(231, 4).Value = 169.0714285714286
$ws.Cells.Item(231, 5).Value = 156.1
$ws.Cells.Item(231, 6).Value = 148
$ws.Cells.Item(231, 7).Value = 64.75
$ws.Cells.Item(231, 8).Value = 134
$ws.Cells.Item(231, 9).Value = 67.83799999999999
$ws.Cells.Item(231, 10).Value = 62
$ws.Cells.Item(231, 11).Value = 125
$ws.Cells.Item(231, 12).Value = 324.19
$ws.Cells.Item(231, 13).Value = 120
$ws.Cells.Item(231, 14).Value = 239
$ws.Cells.Item(231, 15).Value = 198
$ws.Cells.Item(231, 16).Value = 181
$ws.Cells.Item(231, 17).Value = 175
$ws.Cells.Item(231, 18).Value = 99
$ws.Cells.Item(231, 19).Value = 160
$ws.Cells.Item(231, 20).Value = 0.3192307692307693
$ws.Cells.Item(231, 21).Value = 64.48999999999999
$ws.Cells.Item(231, 22).Value = 118.25
$ws.Cells.Item(231, 23).Value = 64.48999999999999

# Row 232
$ws.Cells.Item(232, 1).Value = 230
$ws.Cells.Item(232, 2).Value = "03/09/2024"
$ws.Cells.Item(232, 3).Value = 199.1529411764706
$ws.Cells.Item(232, 4).Value = 165.75
$ws.Cells.Item(232, 5).Value = 156.1
$ws.Cells.Item(232, 6).Value = 148
$ws.Cells.Item(232, 7).Value = 64.75
$ws.Cells.Item(232, 8).Value = 134
$ws.Cells.Item(232, 9).Value = 67.83799999999999
$ws.Cells.Item(232, 10).Value = 62
$ws.Cells.Item(232, 11).Value = 125
$ws.Cells.Item(232, 12).Value = 295.6018181818181
$ws.Cells.Item(232, 13).Value = 120
$ws.Cells.Item(232, 14).Value = 233.5
$ws.Cells.Item(232, 15).Value = 200
$ws.Cells.Item(232, 16).Value = 181
$ws.Cells.Item(232, 17).Value = 175
$ws.Cells.Item(232, 18).Value = 99
$ws.Cells.Item(232, 19).Value = 160
$ws.Cells.Item(232, 20).Value = 0.3192307692307693
$ws.Cells.Item(232, 21).Value = 64.48999999999999
$ws.Cells.Item(232, 22).Value = 118.25
$ws.Cells.Item(232, 23).Value = 64.48999999999999

# Row 233
$ws.Cells.Item(233, 1).Value = 231
$ws.Cells.Item(233, 2).Value = "04/09/2024"
$ws.Cells.Item(233, 3).Value = 198.3666666666667
$ws.Cells.Item(233, 4).Value = 166
$ws.Cells.Item(233, 5).Value = 156.5
$ws.Cells.Item(233, 6).Value = 145.5
$ws.Cells.Item(233, 7).Value = 64.75
$ws.Cells.Item(233, 8).Value = 134
$ws.Cells.Item(233, 9).Value = 67.83799999999999
$ws.Cells.Item(233, 10).Value = 62
$ws.Cells.Item(233, 11).Value = 125
$ws.Cells.Item(233, 12).Value = 298.6857142857143
$ws.Cells.Item(233, 13).Value = 120
$ws.Cells.Item(233, 14).Value = 233.5
$ws.Cells.Item(233, 15).Value = 200
$ws.Cells.Item(233, 16).Value = 181
$ws.Cells.Item(233, 17).Value = 175
$ws.Cells.Item(233, 18).Value = 99
$ws.Cells.Item(233, 19).Value = 160
$ws.Cells.Item(233, 20).Value = 0.3192307692307693
$ws.Cells.Item(233, 21).Value = 64.48999999999999
$ws.Cells.Item(233, 22).Value = 118.25
$ws.Cells.Item(233, 23).Value = 64.48999999999999

# Drop the temporary Text format now that the literal strings are committed.
$ws.Range("B230:B233").ClearFormats()

# Column A uses a bold/centered/thin-bordered style throughout the table;
# copy it from the row immediately above so the new rows match exactly.
$ws.Range("A229").Copy()
$ws.Range("A230:A233").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host "Added rows 230-233 (A1:W233)"